# Anonymize the two participant placeholder names used throughout the
# transcript: "apple127" -> "id76" and "pear128" -> "id77".
# This covers both the bare-name cells (speaker/description columns)
# and the bracketed in-sentence mentions like "[apple127]" / "[pear128]".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("apple127", "id76")
$ws.Cells.Replace("pear128", "id77")
